# Updates NATMI ligand/receptor metrics on Sheet1 (Ly86-Cd180) with
# recomputed TPM-based values (see commit "update scripts wuth new tpm").
# Columns G-J (ligand avg/total expr + specificity) depend only on the
# "Sending cluster" (col A); columns K-P (receptor side) depend only on
# the "Target cluster" (col D); columns Q-T (edge weights/specificity)
# are simple products of the ligand- and receptor-side columns. All 288
# affected cells (rows 2-26, cols G-T) are written explicitly below with
# their new literal values, using a [double] cast so values in
# scientific notation parse correctly as numbers rather than tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = [double]"1.144722"  # G2 (Ligand average expression value)
$ws.Cells.Item(2, 8).Value = [double]"3.434166"  # H2 (Ligand total expression value)
$ws.Cells.Item(2, 9).Value = [double]"0.002843321461639425"  # I2 (Ligand derived specificity of average expression value)
$ws.Cells.Item(2, 10).Value = [double]"0.002849394892703057"  # J2 (Ligand derived specificity of total expression value)
$ws.Cells.Item(2, 13).Value = [double]"2.632114"  # M2 (Receptor average expression value)
$ws.Cells.Item(2, 14).Value = [double]"7.896342000000001"  # N2 (Receptor total expression value)
$ws.Cells.Item(2, 15).Value = [double]"0.009157018146333049"  # O2 (Receptor derived specificity of average expression value)
$ws.Cells.Item(2, 16).Value = [double]"0.00917615351935201"  # P2 (Receptor derived specificity of total expression value)
$ws.Cells.Item(2, 17).Value = [double]"3.013038802308"  # Q2 (Edge average expression weight)
$ws.Cells.Item(2, 18).Value = [double]"27.11734922077201"  # R2 (Edge total expression weight)
$ws.Cells.Item(2, 19).Value = [double]"2.603634622009042E-05"  # S2 (Edge average expression derived specificity)
$ws.Cells.Item(2, 20).Value = [double]"2.61464849727008E-05"  # T2 (Edge total expression derived specificity)
# Row 3
$ws.Cells.Item(3, 7).Value = [double]"1.144722"  # G3 (Ligand average expression value)
$ws.Cells.Item(3, 8).Value = [double]"3.434166"  # H3 (Ligand total expression value)
$ws.Cells.Item(3, 9).Value = [double]"0.002843321461639425"  # I3 (Ligand derived specificity of average expression value)
$ws.Cells.Item(3, 10).Value = [double]"0.002849394892703057"  # J3 (Ligand derived specificity of total expression value)
$ws.Cells.Item(3, 15).Value = [double]"0.002532434389203008"  # O3 (Receptor derived specificity of average expression value)
$ws.Cells.Item(3, 16).Value = [double]"0.002537726404126322"  # P3 (Receptor derived specificity of total expression value)
$ws.Cells.Item(3, 17).Value = [double]"0.833275959164"  # Q3 (Edge average expression weight)
$ws.Cells.Item(3, 18).Value = [double]"7.499483632476001"  # R3 (Edge total expression weight)
$ws.Cells.Item(3, 19).Value = [double]"7.200525049014642E-06"  # S3 (Edge average expression derived specificity)
$ws.Cells.Item(3, 20).Value = [double]"7.230984654995236E-06"  # T3 (Edge total expression derived specificity)
# Row 4
$ws.Cells.Item(4, 7).Value = [double]"1.144722"  # G4 (Ligand average expression value)
$ws.Cells.Item(4, 8).Value = [double]"3.434166"  # H4 (Ligand total expression value)
$ws.Cells.Item(4, 9).Value = [double]"0.002843321461639425"  # I4 (Ligand derived specificity of average expression value)
$ws.Cells.Item(4, 10).Value = [double]"0.002849394892703057"  # J4 (Ligand derived specificity of total expression value)
$ws.Cells.Item(4, 13).Value = [double]"153.046158"  # M4 (Receptor average expression value)
$ws.Cells.Item(4, 14).Value = [double]"459.138474"  # N4 (Receptor total expression value)
$ws.Cells.Item(4, 15).Value = [double]"0.5324413935082427"  # O4 (Receptor derived specificity of average expression value)
$ws.Cells.Item(4, 16).Value = [double]"0.5335540335088085"  # P4 (Receptor derived specificity of total expression value)
$ws.Cells.Item(4, 17).Value = [double]"175.195304078076"  # Q4 (Edge average expression weight)
$ws.Cells.Item(4, 18).Value = [double]"1576.757736702684"  # R4 (Edge total expression weight)
$ws.Cells.Item(4, 19).Value = [double]"0.001513902041227189"  # S4 (Edge average expression derived specificity)
$ws.Cells.Item(4, 20).Value = [double]"0.001520306138061115"  # T4 (Edge total expression derived specificity)
# Row 5
$ws.Cells.Item(5, 7).Value = [double]"1.144722"  # G5 (Ligand average expression value)
$ws.Cells.Item(5, 8).Value = [double]"3.434166"  # H5 (Ligand total expression value)
$ws.Cells.Item(5, 9).Value = [double]"0.002843321461639425"  # I5 (Ligand derived specificity of average expression value)
$ws.Cells.Item(5, 10).Value = [double]"0.002849394892703057"  # J5 (Ligand derived specificity of total expression value)
$ws.Cells.Item(5, 11).Value = [double]"1"  # K5 (Receptor-expressing cells)
$ws.Cells.Item(5, 12).Value = [double]"0.5"  # L5 (Receptor detection rate)
$ws.Cells.Item(5, 13).Value = [double]"1.798242"  # M5 (Receptor average expression value)
$ws.Cells.Item(5, 14).Value = [double]"3.596484"  # N5 (Receptor total expression value)
$ws.Cells.Item(5, 15).Value = [double]"0.006256011185495094"  # O5 (Receptor derived specificity of average expression value)
$ws.Cells.Item(5, 16).Value = [double]"0.004179389559607879"  # P5 (Receptor derived specificity of total expression value)
$ws.Cells.Item(5, 17).Value = [double]"2.058487178724"  # Q5 (Edge average expression weight)
$ws.Cells.Item(5, 18).Value = [double]"12.350923072344"  # R5 (Edge total expression weight)
$ws.Cells.Item(5, 19).Value = [double]"1.77878508679745E-05"  # S5 (Edge average expression derived specificity)
$ws.Cells.Item(5, 20).Value = [double]"1.190873126576317E-05"  # T5 (Edge total expression derived specificity)
# Row 6
$ws.Cells.Item(6, 7).Value = [double]"1.144722"  # G6 (Ligand average expression value)
$ws.Cells.Item(6, 8).Value = [double]"3.434166"  # H6 (Ligand total expression value)
$ws.Cells.Item(6, 9).Value = [double]"0.002843321461639425"  # I6 (Ligand derived specificity of average expression value)
$ws.Cells.Item(6, 10).Value = [double]"0.002849394892703057"  # J6 (Ligand derived specificity of total expression value)
$ws.Cells.Item(6, 13).Value = [double]"129.2378183333333"  # M6 (Receptor average expression value)
$ws.Cells.Item(6, 14).Value = [double]"387.713455"  # N6 (Receptor total expression value)
$ws.Cells.Item(6, 15).Value = [double]"0.4496131427707262"  # O6 (Receptor derived specificity of average expression value)
$ws.Cells.Item(6, 16).Value = [double]"0.4505526970081055"  # P6 (Receptor derived specificity of total expression value)
$ws.Cells.Item(6, 17).Value = [double]"147.94137387817"  # Q6 (Edge average expression weight)
$ws.Cells.Item(6, 18).Value = [double]"1331.47236490353"  # R6 (Edge total expression weight)
$ws.Cells.Item(6, 19).Value = [double]"0.001278394698275157"  # S6 (Edge average expression derived specificity)
$ws.Cells.Item(6, 20).Value = [double]"0.001283802553748484"  # T6 (Edge total expression derived specificity)
# Row 7
$ws.Cells.Item(7, 9).Value = [double]"0.001907216783667987"  # I7 (Ligand derived specificity of average expression value)
$ws.Cells.Item(7, 10).Value = [double]"0.001911290663394668"  # J7 (Ligand derived specificity of total expression value)
$ws.Cells.Item(7, 13).Value = [double]"2.632114"  # M7 (Receptor average expression value)
$ws.Cells.Item(7, 14).Value = [double]"7.896342000000001"  # N7 (Receptor total expression value)
$ws.Cells.Item(7, 15).Value = [double]"0.009157018146333049"  # O7 (Receptor derived specificity of average expression value)
$ws.Cells.Item(7, 16).Value = [double]"0.00917615351935201"  # P7 (Receptor derived specificity of total expression value)
$ws.Cells.Item(7, 17).Value = [double]"2.021058206444"  # Q7 (Edge average expression weight)
$ws.Cells.Item(7, 18).Value = [double]"18.189523857996"  # R7 (Edge total expression weight)
$ws.Cells.Item(7, 19).Value = [double]"1.746441869703871E-05"  # S7 (Edge average expression derived specificity)
$ws.Cells.Item(7, 20).Value = [double]"1.753829654741362E-05"  # T7 (Edge total expression derived specificity)
# Row 8
$ws.Cells.Item(8, 9).Value = [double]"0.001907216783667987"  # I8 (Ligand derived specificity of average expression value)
$ws.Cells.Item(8, 10).Value = [double]"0.001911290663394668"  # J8 (Ligand derived specificity of total expression value)
$ws.Cells.Item(8, 15).Value = [double]"0.002532434389203008"  # O8 (Receptor derived specificity of average expression value)
$ws.Cells.Item(8, 16).Value = [double]"0.002537726404126322"  # P8 (Receptor derived specificity of total expression value)
$ws.Cells.Item(8, 19).Value = [double]"4.829901370625966E-06"  # S8 (Edge average expression derived specificity)
$ws.Cells.Item(8, 20).Value = [double]"4.850332782456764E-06"  # T8 (Edge total expression derived specificity)
# Row 9
$ws.Cells.Item(9, 9).Value = [double]"0.001907216783667987"  # I9 (Ligand derived specificity of average expression value)
$ws.Cells.Item(9, 10).Value = [double]"0.001911290663394668"  # J9 (Ligand derived specificity of total expression value)
$ws.Cells.Item(9, 13).Value = [double]"153.046158"  # M9 (Receptor average expression value)
$ws.Cells.Item(9, 14).Value = [double]"459.138474"  # N9 (Receptor total expression value)
$ws.Cells.Item(9, 15).Value = [double]"0.5324413935082427"  # O9 (Receptor derived specificity of average expression value)
$ws.Cells.Item(9, 16).Value = [double]"0.5335540335088085"  # P9 (Receptor derived specificity of total expression value)
$ws.Cells.Item(9, 17).Value = [double]"117.515880235668"  # Q9 (Edge average expression weight)
$ws.Cells.Item(9, 18).Value = [double]"1057.642922121012"  # R9 (Edge total expression weight)
$ws.Cells.Item(9, 19).Value = [double]"0.001015481162018492"  # S9 (Edge average expression derived specificity)
$ws.Cells.Item(9, 20).Value = [double]"0.001019776842661952"  # T9 (Edge total expression derived specificity)
# Row 10
$ws.Cells.Item(10, 9).Value = [double]"0.001907216783667987"  # I10 (Ligand derived specificity of average expression value)
$ws.Cells.Item(10, 10).Value = [double]"0.001911290663394668"  # J10 (Ligand derived specificity of total expression value)
$ws.Cells.Item(10, 11).Value = [double]"1"  # K10 (Receptor-expressing cells)
$ws.Cells.Item(10, 12).Value = [double]"0.5"  # L10 (Receptor detection rate)
$ws.Cells.Item(10, 13).Value = [double]"1.798242"  # M10 (Receptor average expression value)
$ws.Cells.Item(10, 14).Value = [double]"3.596484"  # N10 (Receptor total expression value)
$ws.Cells.Item(10, 15).Value = [double]"0.006256011185495094"  # O10 (Receptor derived specificity of average expression value)
$ws.Cells.Item(10, 16).Value = [double]"0.004179389559607879"  # P10 (Receptor derived specificity of total expression value)
$ws.Cells.Item(10, 17).Value = [double]"1.380772926732"  # Q10 (Edge average expression weight)
$ws.Cells.Item(10, 18).Value = [double]"8.284637560392"  # R10 (Edge total expression weight)
$ws.Cells.Item(10, 19).Value = [double]"1.19315695317909E-05"  # S10 (Edge average expression derived specificity)
$ws.Cells.Item(10, 20).Value = [double]"7.988028243967692E-06"  # T10 (Edge total expression derived specificity)
# Row 11
$ws.Cells.Item(11, 9).Value = [double]"0.001907216783667987"  # I11 (Ligand derived specificity of average expression value)
$ws.Cells.Item(11, 10).Value = [double]"0.001911290663394668"  # J11 (Ligand derived specificity of total expression value)
$ws.Cells.Item(11, 13).Value = [double]"129.2378183333333"  # M11 (Receptor average expression value)
$ws.Cells.Item(11, 14).Value = [double]"387.713455"  # N11 (Receptor total expression value)
$ws.Cells.Item(11, 15).Value = [double]"0.4496131427707262"  # O11 (Receptor derived specificity of average expression value)
$ws.Cells.Item(11, 16).Value = [double]"0.4505526970081055"  # P11 (Receptor derived specificity of total expression value)
$ws.Cells.Item(11, 17).Value = [double]"99.23474185597668"  # Q11 (Edge average expression weight)
$ws.Cells.Item(11, 18).Value = [double]"893.11267670379"  # R11 (Edge total expression weight)
$ws.Cells.Item(11, 19).Value = [double]"0.0008575097320500401"  # S11 (Edge average expression derived specificity)
$ws.Cells.Item(11, 20).Value = [double]"0.0008611371631588788"  # T11 (Edge total expression derived specificity)
# Row 12
$ws.Cells.Item(12, 7).Value = [double]"211.0125426666667"  # G12 (Ligand average expression value)
$ws.Cells.Item(12, 8).Value = [double]"633.037628"  # H12 (Ligand total expression value)
$ws.Cells.Item(12, 9).Value = [double]"0.5241241901869959"  # I12 (Ligand derived specificity of average expression value)
$ws.Cells.Item(12, 10).Value = [double]"0.5252437372311233"  # J12 (Ligand derived specificity of total expression value)
$ws.Cells.Item(12, 13).Value = [double]"2.632114"  # M12 (Receptor average expression value)
$ws.Cells.Item(12, 14).Value = [double]"7.896342000000001"  # N12 (Receptor total expression value)
$ws.Cells.Item(12, 15).Value = [double]"0.009157018146333049"  # O12 (Receptor derived specificity of average expression value)
$ws.Cells.Item(12, 16).Value = [double]"0.00917615351935201"  # P12 (Receptor derived specificity of total expression value)
$ws.Cells.Item(12, 17).Value = [double]"555.4090677285308"  # Q12 (Edge average expression weight)
$ws.Cells.Item(12, 18).Value = [double]"4998.681609556776"  # R12 (Edge total expression weight)
$ws.Cells.Item(12, 19).Value = [double]"0.004799414720474435"  # S12 (Edge average expression derived specificity)
$ws.Cells.Item(12, 20).Value = [double]"0.004819717167910974"  # T12 (Edge total expression derived specificity)
# Row 13
$ws.Cells.Item(13, 7).Value = [double]"211.0125426666667"  # G13 (Ligand average expression value)
$ws.Cells.Item(13, 8).Value = [double]"633.037628"  # H13 (Ligand total expression value)
$ws.Cells.Item(13, 9).Value = [double]"0.5241241901869959"  # I13 (Ligand derived specificity of average expression value)
$ws.Cells.Item(13, 10).Value = [double]"0.5252437372311233"  # J13 (Ligand derived specificity of total expression value)
$ws.Cells.Item(13, 15).Value = [double]"0.002532434389203008"  # O13 (Receptor derived specificity of average expression value)
$ws.Cells.Item(13, 16).Value = [double]"0.002537726404126322"  # P13 (Receptor derived specificity of total expression value)
$ws.Cells.Item(13, 17).Value = [double]"153.6020788332898"  # Q13 (Edge average expression weight)
$ws.Cells.Item(13, 18).Value = [double]"1382.418709499608"  # R13 (Edge total expression weight)
$ws.Cells.Item(13, 19).Value = [double]"0.001327310123442726"  # S13 (Edge average expression derived specificity)
$ws.Cells.Item(13, 20).Value = [double]"0.001332924900573409"  # T13 (Edge total expression derived specificity)
# Row 14
$ws.Cells.Item(14, 7).Value = [double]"211.0125426666667"  # G14 (Ligand average expression value)
$ws.Cells.Item(14, 8).Value = [double]"633.037628"  # H14 (Ligand total expression value)
$ws.Cells.Item(14, 9).Value = [double]"0.5241241901869959"  # I14 (Ligand derived specificity of average expression value)
$ws.Cells.Item(14, 10).Value = [double]"0.5252437372311233"  # J14 (Ligand derived specificity of total expression value)
$ws.Cells.Item(14, 13).Value = [double]"153.046158"  # M14 (Receptor average expression value)
$ws.Cells.Item(14, 14).Value = [double]"459.138474"  # N14 (Receptor total expression value)
$ws.Cells.Item(14, 15).Value = [double]"0.5324413935082427"  # O14 (Receptor derived specificity of average expression value)
$ws.Cells.Item(14, 16).Value = [double]"0.5335540335088085"  # P14 (Receptor derived specificity of total expression value)
$ws.Cells.Item(14, 17).Value = [double]"32294.65894494441"  # Q14 (Edge average expression weight)
$ws.Cells.Item(14, 18).Value = [double]"290651.9305044997"  # R14 (Edge total expression weight)
$ws.Cells.Item(14, 19).Value = [double]"0.2790654141945433"  # S14 (Edge average expression derived specificity)
$ws.Cells.Item(14, 20).Value = [double]"0.2802459145749066"  # T14 (Edge total expression derived specificity)
# Row 15
$ws.Cells.Item(15, 7).Value = [double]"211.0125426666667"  # G15 (Ligand average expression value)
$ws.Cells.Item(15, 8).Value = [double]"633.037628"  # H15 (Ligand total expression value)
$ws.Cells.Item(15, 9).Value = [double]"0.5241241901869959"  # I15 (Ligand derived specificity of average expression value)
$ws.Cells.Item(15, 10).Value = [double]"0.5252437372311233"  # J15 (Ligand derived specificity of total expression value)
$ws.Cells.Item(15, 11).Value = [double]"1"  # K15 (Receptor-expressing cells)
$ws.Cells.Item(15, 12).Value = [double]"0.5"  # L15 (Receptor detection rate)
$ws.Cells.Item(15, 13).Value = [double]"1.798242"  # M15 (Receptor average expression value)
$ws.Cells.Item(15, 14).Value = [double]"3.596484"  # N15 (Receptor total expression value)
$ws.Cells.Item(15, 15).Value = [double]"0.006256011185495094"  # O15 (Receptor derived specificity of average expression value)
$ws.Cells.Item(15, 16).Value = [double]"0.004179389559607879"  # P15 (Receptor derived specificity of total expression value)
$ws.Cells.Item(15, 17).Value = [double]"379.451616749992"  # Q15 (Edge average expression weight)
$ws.Cells.Item(15, 18).Value = [double]"2276.709700499952"  # R15 (Edge total expression weight)
$ws.Cells.Item(15, 19).Value = [double]"0.003278926796398404"  # S15 (Edge average expression derived specificity)
$ws.Cells.Item(15, 20).Value = [double]"0.00219519819163318"  # T15 (Edge total expression derived specificity)
# Row 16
$ws.Cells.Item(16, 7).Value = [double]"211.0125426666667"  # G16 (Ligand average expression value)
$ws.Cells.Item(16, 8).Value = [double]"633.037628"  # H16 (Ligand total expression value)
$ws.Cells.Item(16, 9).Value = [double]"0.5241241901869959"  # I16 (Ligand derived specificity of average expression value)
$ws.Cells.Item(16, 10).Value = [double]"0.5252437372311233"  # J16 (Ligand derived specificity of total expression value)
$ws.Cells.Item(16, 13).Value = [double]"129.2378183333333"  # M16 (Receptor average expression value)
$ws.Cells.Item(16, 14).Value = [double]"387.713455"  # N16 (Receptor total expression value)
$ws.Cells.Item(16, 15).Value = [double]"0.4496131427707262"  # O16 (Receptor derived specificity of average expression value)
$ws.Cells.Item(16, 16).Value = [double]"0.4505526970081055"  # P16 (Receptor derived specificity of total expression value)
$ws.Cells.Item(16, 17).Value = [double]"27270.80065520942"  # Q16 (Edge average expression weight)
$ws.Cells.Item(16, 18).Value = [double]"245437.2058968848"  # R16 (Edge total expression weight)
$ws.Cells.Item(16, 19).Value = [double]"0.235653124352137"  # S16 (Edge average expression derived specificity)
$ws.Cells.Item(16, 20).Value = [double]"0.2366499823960992"  # T16 (Edge total expression derived specificity)
# Row 17
$ws.Cells.Item(17, 7).Value = [double]"2.5744045"  # G17 (Ligand average expression value)
$ws.Cells.Item(17, 8).Value = [double]"5.148809"  # H17 (Ligand total expression value)
$ws.Cells.Item(17, 9).Value = [double]"0.006394442987722008"  # I17 (Ligand derived specificity of average expression value)
$ws.Cells.Item(17, 10).Value = [double]"0.004272067823193035"  # J17 (Ligand derived specificity of total expression value)
$ws.Cells.Item(17, 13).Value = [double]"2.632114"  # M17 (Receptor average expression value)
$ws.Cells.Item(17, 14).Value = [double]"7.896342000000001"  # N17 (Receptor total expression value)
$ws.Cells.Item(17, 15).Value = [double]"0.009157018146333049"  # O17 (Receptor derived specificity of average expression value)
$ws.Cells.Item(17, 16).Value = [double]"0.00917615351935201"  # P17 (Receptor derived specificity of total expression value)
$ws.Cells.Item(17, 17).Value = [double]"6.776126126113"  # Q17 (Edge average expression weight)
$ws.Cells.Item(17, 18).Value = [double]"40.656756756678"  # R17 (Edge total expression weight)
$ws.Cells.Item(17, 19).Value = [double]"5.855403047426255E-05"  # S17 (Edge average expression derived specificity)
$ws.Cells.Item(17, 20).Value = [double]"3.920115019070326E-05"  # T17 (Edge total expression derived specificity)
# Row 18
$ws.Cells.Item(18, 7).Value = [double]"2.5744045"  # G18 (Ligand average expression value)
$ws.Cells.Item(18, 8).Value = [double]"5.148809"  # H18 (Ligand total expression value)
$ws.Cells.Item(18, 9).Value = [double]"0.006394442987722008"  # I18 (Ligand derived specificity of average expression value)
$ws.Cells.Item(18, 10).Value = [double]"0.004272067823193035"  # J18 (Ligand derived specificity of total expression value)
$ws.Cells.Item(18, 15).Value = [double]"0.002532434389203008"  # O18 (Receptor derived specificity of average expression value)
$ws.Cells.Item(18, 16).Value = [double]"0.002537726404126322"  # P18 (Receptor derived specificity of total expression value)
$ws.Cells.Item(18, 17).Value = [double]"1.873982835145667"  # Q18 (Edge average expression weight)
$ws.Cells.Item(18, 18).Value = [double]"11.243897010874"  # R18 (Edge total expression weight)
$ws.Cells.Item(18, 19).Value = [double]"1.619350732190524E-05"  # S18 (Edge average expression derived specificity)
$ws.Cells.Item(18, 20).Value = [double]"1.084133931513543E-05"  # T18 (Edge total expression derived specificity)
# Row 19
$ws.Cells.Item(19, 7).Value = [double]"2.5744045"  # G19 (Ligand average expression value)
$ws.Cells.Item(19, 8).Value = [double]"5.148809"  # H19 (Ligand total expression value)
$ws.Cells.Item(19, 9).Value = [double]"0.006394442987722008"  # I19 (Ligand derived specificity of average expression value)
$ws.Cells.Item(19, 10).Value = [double]"0.004272067823193035"  # J19 (Ligand derived specificity of total expression value)
$ws.Cells.Item(19, 13).Value = [double]"153.046158"  # M19 (Receptor average expression value)
$ws.Cells.Item(19, 14).Value = [double]"459.138474"  # N19 (Receptor total expression value)
$ws.Cells.Item(19, 15).Value = [double]"0.5324413935082427"  # O19 (Receptor derived specificity of average expression value)
$ws.Cells.Item(19, 16).Value = [double]"0.5335540335088085"  # P19 (Receptor derived specificity of total expression value)
$ws.Cells.Item(19, 17).Value = [double]"394.002717862911"  # Q19 (Edge average expression weight)
$ws.Cells.Item(19, 18).Value = [double]"2364.016307177466"  # R19 (Edge total expression weight)
$ws.Cells.Item(19, 19).Value = [double]"0.003404666135091717"  # S19 (Edge average expression derived specificity)
$ws.Cells.Item(19, 20).Value = [double]"0.00227937901848784"  # T19 (Edge total expression derived specificity)
# Row 20
$ws.Cells.Item(20, 7).Value = [double]"2.5744045"  # G20 (Ligand average expression value)
$ws.Cells.Item(20, 8).Value = [double]"5.148809"  # H20 (Ligand total expression value)
$ws.Cells.Item(20, 9).Value = [double]"0.006394442987722008"  # I20 (Ligand derived specificity of average expression value)
$ws.Cells.Item(20, 10).Value = [double]"0.004272067823193035"  # J20 (Ligand derived specificity of total expression value)
$ws.Cells.Item(20, 11).Value = [double]"1"  # K20 (Receptor-expressing cells)
$ws.Cells.Item(20, 12).Value = [double]"0.5"  # L20 (Receptor detection rate)
$ws.Cells.Item(20, 13).Value = [double]"1.798242"  # M20 (Receptor average expression value)
$ws.Cells.Item(20, 14).Value = [double]"3.596484"  # N20 (Receptor total expression value)
$ws.Cells.Item(20, 15).Value = [double]"0.006256011185495094"  # O20 (Receptor derived specificity of average expression value)
$ws.Cells.Item(20, 16).Value = [double]"0.004179389559607879"  # P20 (Receptor derived specificity of total expression value)
$ws.Cells.Item(20, 17).Value = [double]"4.629402296888999"  # Q20 (Edge average expression weight)
$ws.Cells.Item(20, 18).Value = [double]"18.517609187556"  # R20 (Edge total expression weight)
$ws.Cells.Item(20, 19).Value = [double]"4.000370685619955E-05"  # S20 (Edge average expression derived specificity)
$ws.Cells.Item(20, 20).Value = [double]"1.785463565818973E-05"  # T20 (Edge total expression derived specificity)
# Row 21
$ws.Cells.Item(21, 7).Value = [double]"2.5744045"  # G21 (Ligand average expression value)
$ws.Cells.Item(21, 8).Value = [double]"5.148809"  # H21 (Ligand total expression value)
$ws.Cells.Item(21, 9).Value = [double]"0.006394442987722008"  # I21 (Ligand derived specificity of average expression value)
$ws.Cells.Item(21, 10).Value = [double]"0.004272067823193035"  # J21 (Ligand derived specificity of total expression value)
$ws.Cells.Item(21, 13).Value = [double]"129.2378183333333"  # M21 (Receptor average expression value)
$ws.Cells.Item(21, 14).Value = [double]"387.713455"  # N21 (Receptor total expression value)
$ws.Cells.Item(21, 15).Value = [double]"0.4496131427707262"  # O21 (Receptor derived specificity of average expression value)
$ws.Cells.Item(21, 16).Value = [double]"0.4505526970081055"  # P21 (Receptor derived specificity of total expression value)
$ws.Cells.Item(21, 17).Value = [double]"332.7104210875158"  # Q21 (Edge average expression weight)
$ws.Cells.Item(21, 18).Value = [double]"1996.262526525095"  # R21 (Edge total expression weight)
$ws.Cells.Item(21, 19).Value = [double]"0.002875025607977924"  # S21 (Edge average expression derived specificity)
$ws.Cells.Item(21, 20).Value = [double]"0.001924791679541168"  # T21 (Edge total expression derived specificity)
# Row 22
$ws.Cells.Item(22, 7).Value = [double]"187.100759"  # G22 (Ligand average expression value)
$ws.Cells.Item(22, 8).Value = [double]"561.302277"  # H22 (Ligand total expression value)
$ws.Cells.Item(22, 9).Value = [double]"0.4647308285799747"  # I22 (Ligand derived specificity of average expression value)
$ws.Cells.Item(22, 10).Value = [double]"0.4657235093895858"  # J22 (Ligand derived specificity of total expression value)
$ws.Cells.Item(22, 13).Value = [double]"2.632114"  # M22 (Receptor average expression value)
$ws.Cells.Item(22, 14).Value = [double]"7.896342000000001"  # N22 (Receptor total expression value)
$ws.Cells.Item(22, 15).Value = [double]"0.009157018146333049"  # O22 (Receptor derived specificity of average expression value)
$ws.Cells.Item(22, 16).Value = [double]"0.00917615351935201"  # P22 (Receptor derived specificity of total expression value)
$ws.Cells.Item(22, 17).Value = [double]"492.470527174526"  # Q22 (Edge average expression weight)
$ws.Cells.Item(22, 18).Value = [double]"4432.234744570735"  # R22 (Edge total expression weight)
$ws.Cells.Item(22, 19).Value = [double]"0.004255548630467221"  # S22 (Edge average expression derived specificity)
$ws.Cells.Item(22, 20).Value = [double]"0.004273550419730217"  # T22 (Edge total expression derived specificity)
# Row 23
$ws.Cells.Item(23, 7).Value = [double]"187.100759"  # G23 (Ligand average expression value)
$ws.Cells.Item(23, 8).Value = [double]"561.302277"  # H23 (Ligand total expression value)
$ws.Cells.Item(23, 9).Value = [double]"0.4647308285799747"  # I23 (Ligand derived specificity of average expression value)
$ws.Cells.Item(23, 10).Value = [double]"0.4657235093895858"  # J23 (Ligand derived specificity of total expression value)
$ws.Cells.Item(23, 15).Value = [double]"0.002532434389203008"  # O23 (Receptor derived specificity of average expression value)
$ws.Cells.Item(23, 16).Value = [double]"0.002537726404126322"  # P23 (Receptor derived specificity of total expression value)
$ws.Cells.Item(23, 17).Value = [double]"136.1960060311913"  # Q23 (Edge average expression weight)
$ws.Cells.Item(23, 18).Value = [double]"1225.764054280722"  # R23 (Edge total expression weight)
$ws.Cells.Item(23, 19).Value = [double]"0.001176900332018736"  # S23 (Edge average expression derived specificity)
$ws.Cells.Item(23, 20).Value = [double]"0.001181878846800325"  # T23 (Edge total expression derived specificity)
# Row 24
$ws.Cells.Item(24, 7).Value = [double]"187.100759"  # G24 (Ligand average expression value)
$ws.Cells.Item(24, 8).Value = [double]"561.302277"  # H24 (Ligand total expression value)
$ws.Cells.Item(24, 9).Value = [double]"0.4647308285799747"  # I24 (Ligand derived specificity of average expression value)
$ws.Cells.Item(24, 10).Value = [double]"0.4657235093895858"  # J24 (Ligand derived specificity of total expression value)
$ws.Cells.Item(24, 13).Value = [double]"153.046158"  # M24 (Receptor average expression value)
$ws.Cells.Item(24, 14).Value = [double]"459.138474"  # N24 (Receptor total expression value)
$ws.Cells.Item(24, 15).Value = [double]"0.5324413935082427"  # O24 (Receptor derived specificity of average expression value)
$ws.Cells.Item(24, 16).Value = [double]"0.5335540335088085"  # P24 (Receptor derived specificity of total expression value)
$ws.Cells.Item(24, 17).Value = [double]"28635.05232383392"  # Q24 (Edge average expression weight)
$ws.Cells.Item(24, 18).Value = [double]"257715.4709145053"  # R24 (Edge total expression weight)
$ws.Cells.Item(24, 19).Value = [double]"0.247441929975362"  # S24 (Edge average expression derived specificity)
$ws.Cells.Item(24, 20).Value = [double]"0.248488656934691"  # T24 (Edge total expression derived specificity)
# Row 25
$ws.Cells.Item(25, 7).Value = [double]"187.100759"  # G25 (Ligand average expression value)
$ws.Cells.Item(25, 8).Value = [double]"561.302277"  # H25 (Ligand total expression value)
$ws.Cells.Item(25, 9).Value = [double]"0.4647308285799747"  # I25 (Ligand derived specificity of average expression value)
$ws.Cells.Item(25, 10).Value = [double]"0.4657235093895858"  # J25 (Ligand derived specificity of total expression value)
$ws.Cells.Item(25, 11).Value = [double]"1"  # K25 (Receptor-expressing cells)
$ws.Cells.Item(25, 12).Value = [double]"0.5"  # L25 (Receptor detection rate)
$ws.Cells.Item(25, 13).Value = [double]"1.798242"  # M25 (Receptor average expression value)
$ws.Cells.Item(25, 14).Value = [double]"3.596484"  # N25 (Receptor total expression value)
$ws.Cells.Item(25, 15).Value = [double]"0.006256011185495094"  # O25 (Receptor derived specificity of average expression value)
$ws.Cells.Item(25, 16).Value = [double]"0.004179389559607879"  # P25 (Receptor derived specificity of total expression value)
$ws.Cells.Item(25, 17).Value = [double]"336.452443065678"  # Q25 (Edge average expression weight)
$ws.Cells.Item(25, 18).Value = [double]"2018.714658394068"  # R25 (Edge total expression weight)
$ws.Cells.Item(25, 19).Value = [double]"0.002907361261840725"  # S25 (Edge average expression derived specificity)
$ws.Cells.Item(25, 20).Value = [double]"0.001946439972806777"  # T25 (Edge total expression derived specificity)
# Row 26
$ws.Cells.Item(26, 7).Value = [double]"187.100759"  # G26 (Ligand average expression value)
$ws.Cells.Item(26, 8).Value = [double]"561.302277"  # H26 (Ligand total expression value)
$ws.Cells.Item(26, 9).Value = [double]"0.4647308285799747"  # I26 (Ligand derived specificity of average expression value)
$ws.Cells.Item(26, 10).Value = [double]"0.4657235093895858"  # J26 (Ligand derived specificity of total expression value)
$ws.Cells.Item(26, 13).Value = [double]"129.2378183333333"  # M26 (Receptor average expression value)
$ws.Cells.Item(26, 14).Value = [double]"387.713455"  # N26 (Receptor total expression value)
$ws.Cells.Item(26, 15).Value = [double]"0.4496131427707262"  # O26 (Receptor derived specificity of average expression value)
$ws.Cells.Item(26, 16).Value = [double]"0.4505526970081055"  # P26 (Receptor derived specificity of total expression value)
$ws.Cells.Item(26, 17).Value = [double]"24180.49390167078"  # Q26 (Edge average expression weight)
$ws.Cells.Item(26, 18).Value = [double]"217624.445115037"  # R26 (Edge total expression weight)
$ws.Cells.Item(26, 19).Value = [double]"0.2089490883802861"  # S26 (Edge average expression derived specificity)
$ws.Cells.Item(26, 20).Value = [double]"0.2098329832155576"  # T26 (Edge total expression derived specificity)
